$wb = $excel.ActiveWorkbook

$wsSuccess = $wb.Worksheets.Item("ISPUserSuccessfulLogin")
$wsInvalid = $wb.Worksheets.Item("ISPLoginWithInvalidPswd")

# Remove the duplicate/stale "fadmin" credential row from the successful-login
# sheet (row 2: fadmin/123qwe) - the "admin" row shifts up to row 2.
$wsSuccess.Rows.Item(2).Delete()

# Refresh selections to match the refactored sheets.
$wsInvalid.Activate()
$wsInvalid.Range("B2").Select()

$wsSuccess.Activate()
$wsSuccess.Range("E13").Select()
